$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 45; this shifts the existing rows 45-55 down to 46-56
# and extends the used range to row 56 (matches the updated <dimension ref="A1:R56"/>).
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly record.
$ws.Range("A45").Value = 5
$ws.Range("B45").Value = "Macroferia Regional de Talca"
$ws.Range("C45").Value = "Maule"
$ws.Range("D45").Value = 44505
$ws.Range("E45").Value = 7
$ws.Range("F45").Value = 100112022
$ws.Range("G45").Value = "Arveja Verde"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 600
$ws.Range("K45").Value = 12000
$ws.Range("L45").Value = 12000
$ws.Range("M45").Value = 12000
$ws.Range("N45").Value = "`$/saco 25 kilos"
$ws.Range("O45").Value = "Región del Maule"
$ws.Range("P45").Value = 480
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"
